$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("BF$r").Value = "2008-04-18"
}
